$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws15 = $wb.Worksheets.Item(15)

# --- Sheet1: collapse the 6-step login sequence (rows 2-7) into a single "login" step,
# then shift remaining rows (old 8-10 -> new 4-6 happened via overwrite below),
# and drop the now-unused trailing rows.
$ws1.Range("A2").Value = 1
$ws1.Range("B2").Value = 'login'
$ws1.Range("C2").Value = 'makes a login'
$ws1.Range("D2").ClearContents()

$ws1.Range("A3").Value = 2
$ws1.Range("B3").Value = 'click'
$ws1.Range("C3").Value = 'Click on Login'
$ws1.Range("E3").Value = 'xpath'
$ws1.Range("F3").Value = '//button[@class=''mat-raised-button'']'

$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = 'click'
$ws1.Range("C4").Value = 'Click on ''Menu'''
$ws1.Range("E4").Value = 'xpath'
$ws1.Range("F4").Value = '/html[1]/body[1]/app-root[1]/app-apphome[1]/app-layout[1]/div[1]/div[1]/mat-sidenav-container[1]/mat-sidenav-content[1]/app-app-header[1]/mat-toolbar[1]/div[1]/button[1]/span[1]/mat-icon[1]'

$ws1.Range("A5").Value = 4
$ws1.Range("B5").Value = 'click'
$ws1.Range("C5").Value = 'Click on ''New Project'''
$ws1.Range("E5").Value = 'xpath'
$ws1.Range("F5").Value = '//span[contains(text(),''New Projct'')]'
$ws1.Range("D5").ClearContents()

$ws1.Range("A6").Value = 5
$ws1.Range("B6").Value = 'quit'
$ws1.Range("C6").ClearContents()
$ws1.Range("D6").ClearContents()
$ws1.Range("E6").ClearContents()
$ws1.Range("F6").ClearContents()

$ws1.Rows("7:10").Delete()

# --- Sheet2: same collapse (rows 2-7 -> a single "login" row), remaining rows overwritten
# in place, trailing now-unused rows removed.
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = 'login'
$ws2.Range("C2").Value = 'makes a login'
$ws2.Range("D2").ClearContents()

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = 'click'
$ws2.Range("C3").Value = 'Click on ''Menu'''
$ws2.Range("E3").Value = 'xpath'
$ws2.Range("F3").Value = '/html[1]/body[1]/app-root[1]/app-apphome[1]/app-layout[1]/div[1]/div[1]/mat-sidenav-container[1]/mat-sidenav-content[1]/app-app-header[1]/mat-toolbar[1]/div[1]/button[1]/span[1]/mat-icon[1]'

$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = 'click'
$ws2.Range("C4").Value = 'Click on ''New Project'''
$ws2.Range("E4").Value = 'xpath'
$ws2.Range("F4").Value = '//span[contains(text(),''New Project'')]'

$ws2.Range("A5").Value = 4
$ws2.Range("B5").Value = 'type'
$ws2.Range("C5").Value = 'Type the project name'
$ws2.Range("D5").Value = 'Project1'
$ws2.Range("E5").Value = 'name'
$ws2.Range("F5").Value = 'inName'

$ws2.Range("A6").Value = 5
$ws2.Range("B6").Value = 'type'
$ws2.Range("C6").Value = 'Type the description of the project'
$ws2.Range("D6").Value = 'This is the description of the first project in the sprint'
$ws2.Range("E6").Value = 'name'
$ws2.Range("F6").Value = 'inDesc'

$ws2.Range("A7").Value = 6
$ws2.Range("B7").Value = 'click'
$ws2.Range("C7").Value = 'Select Scrum Master'
$ws2.Range("E7").Value = 'xpath'
$ws2.Range("F7").Value = '//option[contains(text(),''Scrum Master'')]'

$ws2.Range("A8").Value = 7
$ws2.Range("B8").Value = 'type'
$ws2.Range("C8").Value = 'Type the username'
$ws2.Range("D8").Value = 'MarisolR'
$ws2.Range("E8").Value = 'id'
$ws2.Range("F8").Value = 'mat-input-5'

$ws2.Range("A9").Value = 8
$ws2.Range("B9").Value = 'click'
$ws2.Range("C9").Value = 'Click on "Add"'
$ws2.Range("E9").Value = 'xpath'
$ws2.Range("F9").Value = '//span[contains(text(),''Add'')]'

$ws2.Range("A10").Value = 9
$ws2.Range("B10").Value = 'calendar'
$ws2.Range("C10").Value = 'Select calendar'
$ws2.Range("D10").ClearContents()
$ws2.Range("E10").ClearContents()
$ws2.Range("F10").ClearContents()

$ws2.Range("A11").Value = 10
$ws2.Range("B11").Value = 'click'
$ws2.Range("C11").Value = 'Select the start date'
$ws2.Range("E11").Value = 'xpath'
$ws2.Range("F11").Value = '//mat-calendar//div[contains(text(),''23'')]'
$ws2.Range("D11").ClearContents()

$ws2.Range("A12").Value = 11
$ws2.Range("B12").Value = 'click'
$ws2.Range("C12").Value = 'Enable end date'
$ws2.Range("E12").Value = 'xpath'
$ws2.Range("F12").Value = '//div[@class=''mat-checkbox-inner-container'']'

$ws2.Range("A13").Value = 12
$ws2.Range("B13").Value = 'click'
$ws2.Range("C13").Value = 'Select calendar'
$ws2.Range("E13").Value = 'xpath'
$ws2.Range("F13").Value = '/html[1]/body[1]/div[1]/div[2]/div[1]/mat-dialog-container[1]/app-newproject[1]/div[1]/form[1]/div[1]/div[1]/div[2]/div[1]/div[2]/div[1]/mat-form-field[1]/div[1]/div[1]/div[2]/mat-datepicker-toggle[1]/button[1]/span[1]/*[local-name()=''svg''][1]'
$ws2.Range("D13").ClearContents()

$ws2.Range("A14").Value = 13
$ws2.Range("B14").Value = 'click'
$ws2.Range("C14").Value = 'Go to the month April'
$ws2.Range("E14").Value = 'xpath'
$ws2.Range("F14").Value = '//button[@class=''mat-calendar-next-button mat-icon-button'']'

$ws2.Range("A15").Value = 14
$ws2.Range("B15").Value = 'click'
$ws2.Range("C15").Value = 'Select the day'
$ws2.Range("E15").Value = 'xpath'
$ws2.Range("F15").Value = '//mat-calendar//div[contains(text(),''7'')]'

$ws2.Range("A16").Value = 15
$ws2.Range("B16").Value = 'click'
$ws2.Range("C16").Value = 'Click the button "Create"'
$ws2.Range("E16").Value = 'xpath'
$ws2.Range("F16").Value = '//span[contains(text(),''Create'')]'

$ws2.Range("A17").Value = 16
$ws2.Range("B17").Value = 'alert'
$ws2.Range("C17").Value = 'Close the alert'
$ws2.Range("E17").ClearContents()
$ws2.Range("F17").ClearContents()

$ws2.Range("A18").Value = 17
$ws2.Range("B18").Value = 'quit'
$ws2.Range("C18").ClearContents()
$ws2.Range("E18").ClearContents()
$ws2.Range("F18").ClearContents()

$ws2.Rows("19:23").Delete()

# --- View state: selections + active sheet/tab bookkeeping.
$ws1.Range("C7").Select()
$ws15.Range("I6").Select()
$ws2.Activate()
$ws2.Range("F16").Select()
